$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 5.922420666666667
$ws.Range("H2").Value2 = 17.767262
$ws.Range("I2").Value2 = 0.5833698282960434
$ws.Range("J2").Value2 = 0.6311054116979437
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.535712666666667
$ws.Range("N2").Value2 = 7.607138
$ws.Range("O2").Value2 = 0.04494879354621957
$ws.Range("P2").Value2 = 0.05070282964779482
$ws.Range("Q2").Value2 = 15.01755710179511
$ws.Range("R2").Value2 = 135.158013916156
$ws.Range("S2").Value2 = 0.02622176997317242
$ws.Range("T2").Value2 = 0.03199883017912225

$ws.Range("G3").Value2 = 5.922420666666667
$ws.Range("H3").Value2 = 17.767262
$ws.Range("I3").Value2 = 0.5833698282960434
$ws.Range("J3").Value2 = 0.6311054116979437
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 32.24261766666667
$ws.Range("N3").Value2 = 96.72785300000001
$ws.Range("O3").Value2 = 0.5715421877013505
$ws.Range("P3").Value2 = 0.6447070965264385
$ws.Range("Q3").Value2 = 190.9543452164984
$ws.Range("R3").Value2 = 1718.589106948486
$ws.Range("S3").Value2 = 0.3334204679032818
$ws.Range("T3").Value2 = 0.4068781375779039

$ws.Range("G4").Value2 = 5.922420666666667
$ws.Range("H4").Value2 = 17.767262
$ws.Range("I4").Value2 = 0.5833698282960434
$ws.Range("J4").Value2 = 0.6311054116979437
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 1.538811333333333
$ws.Range("N4").Value2 = 4.616434
$ws.Range("O4").Value2 = 0.02727742533206951
$ws.Range("P4").Value2 = 0.03076929413956839
$ws.Range("Q4").Value2 = 9.113488042634222
$ws.Range("R4").Value2 = 82.02139238370799
$ws.Range("S4").Value2 = 0.01591282693232754
$ws.Range("T4").Value2 = 0.01941866804560743

$ws.Range("G5").Value2 = 5.922420666666667
$ws.Range("H5").Value2 = 17.767262
$ws.Range("I5").Value2 = 0.5833698282960434
$ws.Range("J5").Value2 = 0.6311054116979437
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 19.206297
$ws.Range("N5").Value2 = 38.412594
$ws.Range("O5").Value2 = 0.3404565075487166
$ws.Range("P5").Value2 = 0.2560262755732715
$ws.Range("Q5").Value2 = 113.747770282938
$ws.Range("R5").Value2 = 682.4866216976279
$ws.Range("S5").Value2 = 0.1986120543509654
$ws.Range("T5").Value2 = 0.1615795680511607

$ws.Range("G6").Value2 = 5.922420666666667
$ws.Range("H6").Value2 = 17.767262
$ws.Range("I6").Value2 = 0.5833698282960434
$ws.Range("J6").Value2 = 0.6311054116979437
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 0.8899256666666666
$ws.Range("N6").Value2 = 2.669777
$ws.Range("O6").Value2 = 0.0157750858716439
$ws.Range("P6").Value2 = 0.01779450411292666
$ws.Range("Q6").Value2 = 5.270514160063778
$ws.Range("R6").Value2 = 47.43462744057399
$ws.Range("S6").Value2 = 0.009202709136296243
$ws.Range("T6").Value2 = 0.01123020784414933

$ws.Range("G7").Value2 = 1.926013333333334
$ws.Range("H7").Value2 = 5.778040000000001
$ws.Range("I7").Value2 = 0.1897160182974547
$ws.Range("J7").Value2 = 0.2052399696141807
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 2.535712666666667
$ws.Range("N7").Value2 = 7.607138
$ws.Range("O7").Value2 = 0.04494879354621957
$ws.Range("P7").Value2 = 0.05070282964779482
$ws.Range("Q7").Value2 = 4.883816405502222
$ws.Range("R7").Value2 = 43.95434764952
$ws.Range("S7").Value2 = 0.008527506138863106
$ws.Range("T7").Value2 = 0.01040624721626639

$ws.Range("G8").Value2 = 1.926013333333334
$ws.Range("H8").Value2 = 5.778040000000001
$ws.Range("I8").Value2 = 0.1897160182974547
$ws.Range("J8").Value2 = 0.2052399696141807
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 32.24261766666667
$ws.Range("N8").Value2 = 96.72785300000001
$ws.Range("O8").Value2 = 0.5715421877013505
$ws.Range("P8").Value2 = 0.6447070965264385
$ws.Range("Q8").Value2 = 62.0997115275689
$ws.Range("R8").Value2 = 558.8974037481202
$ws.Range("S8").Value2 = 0.1084307081397167
$ws.Range("T8").Value2 = 0.1323196649011329

$ws.Range("G9").Value2 = 1.926013333333334
$ws.Range("H9").Value2 = 5.778040000000001
$ws.Range("I9").Value2 = 0.1897160182974547
$ws.Range("J9").Value2 = 0.2052399696141807
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 1.538811333333333
$ws.Range("N9").Value2 = 4.616434
$ws.Range("O9").Value2 = 0.02727742533206951
$ws.Range("P9").Value2 = 0.03076929413956839
$ws.Range("Q9").Value2 = 2.963771145484445
$ws.Range("R9").Value2 = 26.67394030936
$ws.Range("S9").Value2 = 0.005174964523406354
$ws.Range("T9").Value2 = 0.006315088994254805

$ws.Range("G10").Value2 = 1.926013333333334
$ws.Range("H10").Value2 = 5.778040000000001
$ws.Range("I10").Value2 = 0.1897160182974547
$ws.Range("J10").Value2 = 0.2052399696141807
$ws.Range("K10").Value2 = 2
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 19.206297
$ws.Range("N10").Value2 = 38.412594
$ws.Range("O10").Value2 = 0.3404565075487166
$ws.Range("P10").Value2 = 0.2560262755732715
$ws.Range("Q10").Value2 = 36.99158410596
$ws.Range("R10").Value2 = 221.94950463576
$ws.Range("S10").Value2 = 0.06459005301559984
$ws.Range("T10").Value2 = 0.05254682501909009

$ws.Range("G11").Value2 = 1.926013333333334
$ws.Range("H11").Value2 = 5.778040000000001
$ws.Range("I11").Value2 = 0.1897160182974547
$ws.Range("J11").Value2 = 0.2052399696141807
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 0.8899256666666666
$ws.Range("N11").Value2 = 2.669777
$ws.Range("O11").Value2 = 0.0157750858716439
$ws.Range("P11").Value2 = 0.01779450411292666
$ws.Range("Q11").Value2 = 1.714008699675556
$ws.Range("R11").Value2 = 15.42607829708
$ws.Range("S11").Value2 = 0.002992786479868713
$ws.Range("T11").Value2 = 0.003652143483436481

$ws.Range("G12").Value2 = 2.303652
$ws.Range("H12").Value2 = 4.607303999999999
$ws.Range("I12").Value2 = 0.2269141534065018
$ws.Range("J12").Value2 = 0.1636546186878756
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 2.535712666666667
$ws.Range("N12").Value2 = 7.607138
$ws.Range("O12").Value2 = 0.04494879354621957
$ws.Range("P12").Value2 = 0.05070282964779482
$ws.Range("Q12").Value2 = 5.841399555991998
$ws.Range("R12").Value2 = 35.048397335952
$ws.Range("S12").Value2 = 0.01019951743418405
$ws.Range("T12").Value2 = 0.008297752252406176

$ws.Range("G13").Value2 = 2.303652
$ws.Range("H13").Value2 = 4.607303999999999
$ws.Range("I13").Value2 = 0.2269141534065018
$ws.Range("J13").Value2 = 0.1636546186878756
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 32.24261766666667
$ws.Range("N13").Value2 = 96.72785300000001
$ws.Range("O13").Value2 = 0.5715421877013505
$ws.Range("P13").Value2 = 0.6447070965264385
$ws.Range("Q13").Value2 = 74.27577067305199
$ws.Range("R13").Value2 = 445.654624038312
$ws.Range("S13").Value2 = 0.1296910116583519
$ws.Range("T13").Value2 = 0.1055092940474017

$ws.Range("G14").Value2 = 2.303652
$ws.Range("H14").Value2 = 4.607303999999999
$ws.Range("I14").Value2 = 0.2269141534065018
$ws.Range("J14").Value2 = 0.1636546186878756
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 1.538811333333333
$ws.Range("N14").Value2 = 4.616434
$ws.Range("O14").Value2 = 0.02727742533206951
$ws.Range("P14").Value2 = 0.03076929413956839
$ws.Range("Q14").Value2 = 3.544885805655999
$ws.Range("R14").Value2 = 21.269314833936
$ws.Range("S14").Value2 = 0.006189633876335619
$ws.Range("T14").Value2 = 0.005035537099706151

$ws.Range("G15").Value2 = 2.303652
$ws.Range("H15").Value2 = 4.607303999999999
$ws.Range("I15").Value2 = 0.2269141534065018
$ws.Range("J15").Value2 = 0.1636546186878756
$ws.Range("K15").Value2 = 2
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 19.206297
$ws.Range("N15").Value2 = 38.412594
$ws.Range("O15").Value2 = 0.3404565075487166
$ws.Range("P15").Value2 = 0.2560262755732715
$ws.Range("Q15").Value2 = 44.24462449664399
$ws.Range("R15").Value2 = 176.978497986576
$ws.Range("S15").Value2 = 0.07725440018215132
$ws.Range("T15").Value2 = 0.0418998825030207

$ws.Range("G16").Value2 = 2.303652
$ws.Range("H16").Value2 = 4.607303999999999
$ws.Range("I16").Value2 = 0.2269141534065018
$ws.Range("J16").Value2 = 0.1636546186878756
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 0.8899256666666666
$ws.Range("N16").Value2 = 2.669777
$ws.Range("O16").Value2 = 0.0157750858716439
$ws.Range("P16").Value2 = 0.01779450411292666
$ws.Range("Q16").Value2 = 2.050079041867999
$ws.Range("R16").Value2 = 12.300474251208
$ws.Range("S16").Value2 = 0.009202709136296243
$ws.Range("T16").Value2 = 0.01123020784414933

